$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 46; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C = "Förändrad"
    $cell.Value2 = $cell.Value2 + 1
}
